$wb = $excel.ActiveWorkbook

# --- Settings sheet -------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")

# Update the business-process-name default value.
$wsSettings.Range("B5").Value = "FInanceAndAccounting-REF-CalculateClientSecurityHash"

# New rows: System1 credential + URLs for the new process.
$wsSettings.Range("A6").Value = "System1_CredentialName"
$wsSettings.Range("B6").Value = "ACMELogin"
$wsSettings.Range("C6").Value = "Log in credentials for ACME System 1"

$wsSettings.Range("A7").Value = "System1_URL"
$wsSettings.Range("A8").Value = "SHA1Online_URL"

# Hyperlinks (added in this order so relationship ids line up: rId1 -> B8, rId2 -> B7)
$wsSettings.Hyperlinks.Add($wsSettings.Range("B8"), "http://www.sha1-online.com/")
$wsSettings.Hyperlinks.Add($wsSettings.Range("B7"), "https://acme-test.uipath.com")

# --- Constants sheet --------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("B2").Value = 2

# --- Active sheet / selection bookkeeping -----------------------------------
$null = $wsSettings.Range("A16").Select()
$null = $wsConstants.Activate()
$null = $wsConstants.Range("B3").Select()
